$d = $word.ActiveDocument

# 1) Merge the split "Cam" / "pus Saint-Jean..." runs back into a single run
#    (this also removes the stray _GoBack bookmark that previously sat
#     between them), restoring "Campus Saint-Jean event coordinators."
$d.Content.Find.Execute("Campus Saint-Jean event coordinators.", $true, $false, $false, $false, $false, $true, 1, $false, "Campus Saint-Jean event coordinators.", 2)

# 2) Re-create the _GoBack bookmark at the new "last edit" location: right
#    after the stand-alone "," run in "...basic HTML, minimal CSS, image
#    sprites," and before " and illustrated backgrounds."
$r = $d.Content
$r.Find.Execute("sprites,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
